$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B40 currently holds the text "3"; the update re-stores it as a real number
$ws.Range("B40").Value = 3

# Append a new annotation row (row 41) with the same "Ruilin" annotator data
$ws.Range("A41").Value = "Ruilin"

# B41 must stay a text "3" (unlike B40), so force text storage for the numeric-looking string
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "3"
$ws.Range("B41").ClearFormats()

$ws.Range("C41").Value = "无"
$ws.Range("D41").Value = "DFT"
$ws.Range("E41").Value = "WRI"
$ws.Range("F41").Value = "cf97de89-8b46-4ca2-a071-801296a106cf"
$ws.Range("G41").Value = "SkYXvCR6W_annotated.xlsx"
$ws.Range("H41").Value = "Some important implementation details are missing (activation functions, loss function used), and others have to be deduced by observing the output dimensions of the individual layers of the network."
